$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain decimal numbers must be force-set
# to Text format first, otherwise Excel/COM auto-converts them to numeric values
# (losing e.g. trailing zeros such as "10.10" -> 10.1).

$ws.Range("D2").Value = "26.839.87"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "1.815.08"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.27"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4601"
$ws.Range("E7").Value = "  -2.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3635"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07210"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8565"
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.69"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07508"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "1.763.73"
$ws.Range("E13").Value = "  -8.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.320"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.522"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.63"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008571"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "26.875.88"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.40"
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.137"
$ws.Range("E22").Value = "  -3.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.48"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").Value = "2.033.52"
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.05"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.851"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.12"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.084"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.19"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08846"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.951"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.406"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.130"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7162"
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.072"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05242"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.432"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.918"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.139"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5115"
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1618"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.164"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4786"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.10"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.90"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.614"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06190"
$ws.Range("E51").Value = "  -1.79%  "
